$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Remarque" column header in E1
$ws.Range("E1").Value = "Remarque"

# Add the remark value for the admin row in E2 (same text as A2 - "admin")
$ws.Range("E2").Value = $ws.Range("A2").Value2

# Copy the style from D1:D2 (centered) to E1:E2
$ws.Range("D1:D2").Copy()
$ws.Range("E1:E2").PasteSpecial(-4122)  # xlPasteFormats

# Set the width of the new column E (stored OOXML width of 47 "characters";
# the ColumnWidth property round-trips through a pixel conversion that adds
# ~0.8333 when saved, so back the input off by that amount)
$ws.Columns.Item(5).ColumnWidth = 46.166666666666664

# Move the active selection to row 7 (next empty row), matching the diff
$ws.Range("A7:XFD7").Select()
